$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 3.15
$ws.Range("I2").Value = 4.5
$ws.Range("J2").Value = 2.92
$ws.Range("K2").Value = 6
$ws.Range("T2").Value = 1.66
$ws.Range("U2").Value = 1.72
$ws.Range("V2").Value = 1.29

# Row 3
$ws.Range("I3").Value = 15
$ws.Range("J3").Value = 6.2
$ws.Range("N3").Value = 5
$ws.Range("AH3").Value = 1000

# Row 4
$ws.Range("F4").Value = 1.64
$ws.Range("I4").Value = 12.5
$ws.Range("V4").Value = 1.08
$ws.Range("Y4").Value = 21
$ws.Range("Z4").Value = 980
$ws.Range("AJ4").Value = 28
$ws.Range("AL4").Value = 60
$ws.Range("AN4").Value = 1000

# Row 5
$ws.Range("F5").Value = 2
$ws.Range("G5").Value = 2.18
$ws.Range("J5").Value = 3.2
$ws.Range("K5").Value = 3.6
$ws.Range("N5").Value = 3
$ws.Range("P5").Value = 1.68
$ws.Range("Q5").Value = 2.26
$ws.Range("T5").Value = 1.97
$ws.Range("U5").Value = 1.87
$ws.Range("W5").Value = 1.84
$ws.Range("Z5").Value = 34
$ws.Range("AJ5").Value = 980

# Row 6
$ws.Range("H6").Value = 20
$ws.Range("I6").Value = 27
$ws.Range("J6").Value = 6.6
$ws.Range("K6").Value = 7.8
$ws.Range("M6").Value = 1.03
$ws.Range("Q6").Value = 1.63
$ws.Range("S6").Value = 2.6
$ws.Range("T6").Value = 2.52
$ws.Range("Z6").Value = 270
$ws.Range("AD6").Value = 100
$ws.Range("AI6").Value = 410
$ws.Range("AM6").Value = 440

# Row 7
$ws.Range("F7").Value = 1.31
$ws.Range("G7").Value = 1.36
$ws.Range("H7").Value = 12.5
$ws.Range("I7").Value = 15.5
$ws.Range("J7").Value = 5.3
$ws.Range("K7").Value = 6.6
$ws.Range("L7").Value = 1.36
$ws.Range("N7").Value = 4.2
$ws.Range("O7").Value = 1.27
$ws.Range("Q7").Value = 1.72
$ws.Range("R7").Value = 1.43
$ws.Range("S7").Value = 2.84
$ws.Range("T7").Value = 2.22
$ws.Range("U7").Value = 1.65
$ws.Range("V7").Value = 1.06
$ws.Range("W7").Value = 3.75
$ws.Range("Z7").Value = 150
$ws.Range("AB7").Value = 9.199999999999999
$ws.Range("AC7").Value = 14
$ws.Range("AD7").Value = 980
$ws.Range("AF7").Value = 8.800000000000001
$ws.Range("AH7").Value = 40
$ws.Range("AJ7").Value = 10.5
$ws.Range("AK7").Value = 18.5
$ws.Range("AN7").Value = 6.8
$ws.Range("AO7").Value = 500

# Row 8
$ws.Range("G8").Value = 2
$ws.Range("J8").Value = 3.4
$ws.Range("R8").Value = 1.27
$ws.Range("W8").Value = 2
$ws.Range("AA8").Value = 150
$ws.Range("AD8").Value = 20
$ws.Range("AE8").Value = 80
$ws.Range("AI8").Value = 100
$ws.Range("AK8").Value = 23
$ws.Range("AO8").Value = 120

# Row 9
$ws.Range("F9").Value = 2
$ws.Range("G9").Value = 2.18
$ws.Range("H9").Value = 4.2
$ws.Range("I9").Value = 4.9
$ws.Range("J9").Value = 3.1
$ws.Range("M9").Value = 1.1
$ws.Range("N9").Value = 2.84
$ws.Range("P9").Value = 1.63
$ws.Range("Q9").Value = 2.28
$ws.Range("T9").Value = 1.98
$ws.Range("U9").Value = 1.83
$ws.Range("V9").Value = 1.26
$ws.Range("W9").Value = 1.84
$ws.Range("X9").Value = 10
$ws.Range("Y9").Value = 13.5
$ws.Range("Z9").Value = 32
$ws.Range("AA9").Value = 120
$ws.Range("AB9").Value = 7.8
$ws.Range("AC9").Value = 7.8
$ws.Range("AD9").Value = 20
$ws.Range("AE9").Value = 80
$ws.Range("AF9").Value = 12.5
$ws.Range("AG9").Value = 11.5
$ws.Range("AH9").Value = 23
$ws.Range("AI9").Value = 980
$ws.Range("AJ9").Value = 27
$ws.Range("AK9").Value = 27
$ws.Range("AL9").Value = 55
$ws.Range("AM9").Value = 180
$ws.Range("AN9").Value = 26

# Row 10
$ws.Range("G10").Value = 1.12
$ws.Range("H10").Value = 21
$ws.Range("J10").Value = 13
$ws.Range("K10").Value = 16.5
$ws.Range("N10").Value = 3.7
$ws.Range("O10").Value = 1.08
$ws.Range("P10").Value = 3.7
$ws.Range("Q10").Value = 1.3
$ws.Range("S10").Value = 1.64
$ws.Range("T10").Value = 2.28
$ws.Range("U10").Value = 1.42
$ws.Range("W10").Value = 9.4
$ws.Range("AB10").Value = 1000
$ws.Range("AF10").Value = 12.5
$ws.Range("AG10").Value = 24
$ws.Range("AK10").Value = 22
$ws.Range("AL10").Value = 90

# Row 11
$ws.Range("F11").Value = 2.78
$ws.Range("G11").Value = 3.25
$ws.Range("Q11").Value = 2.08
